# StartCodeAndmemoryMapForGoL.xlsx - "added .org to emulator code" edit
#
# The only real data change is Sheet1!S2: 8 -> 576 (the base/origin address
# used throughout the "memory map" table). Every other changed cell in the
# sheet is a formula (DEC2OCT / "+2" chains) that recalculates automatically
# once S2 changes, so we only need to touch the one input cell.
#
# The rest of the diff unhides the previously-hidden helper columns S:AI
# (and widens column AI to fit), and updates the sheet's selection/scroll
# to showcase that area.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Update the base address input cell; dependent formulas recalc automatically.
$ws.Range("S2").Value = 576

# 2) Unhide the helper columns S:AI that hold the memory-map scratch data.
$ws.Columns("S:AI").Hidden = $false

# Column AI (35) also gets a real width now that it's visible again.
$ws.Columns("AI").ColumnWidth = 8.14

# 3) Update the visible selection to the newly-revealed octal table.
$ws.Range("AL2").Select()
$excel.ActiveWindow.ScrollColumn = 36
$excel.ActiveWindow.ScrollRow = 1
